$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted numbers (e.g. "28.304.88" using
# "." as both thousands + decimal separator) - force text so Excel does not
# auto-coerce plain-looking values (e.g. "210.19") into real numbers.

# --- Update Price (D) and Volume(1h) (E) columns for rows with simple value refresh ---
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "28.304.88"
$ws.Cells.Item(2, 5).Value = "  -1.51%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.551.69"
$ws.Cells.Item(3, 5).Value = "  -1.45%  "
$ws.Cells.Item(4, 5).Value = "  -0.20%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "210.19"
$ws.Cells.Item(5, 5).Value = "  -1.50%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.482"
$ws.Cells.Item(6, 5).Value = "  -1.89%  "
$ws.Cells.Item(7, 5).Value = "  -0.15%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "23.79"
$ws.Cells.Item(8, 5).Value = "  -1.57%  "
$ws.Cells.Item(9, 5).Value = "  -1.84%  "
$ws.Cells.Item(10, 5).Value = "  -1.70%  "
$ws.Cells.Item(11, 5).Value = "  +0.06%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.773.56"
$ws.Cells.Item(12, 5).Value = "  -1.46%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.546.92"
$ws.Cells.Item(13, 5).Value = "  -1.81%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "28.279.64"
$ws.Cells.Item(14, 5).Value = "  -1.61%  "
$ws.Cells.Item(15, 5).Value = "  -1.69%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.510"
$ws.Cells.Item(16, 5).Value = "  -2.54%  "
$ws.Cells.Item(17, 5).Value = "  -3.06%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "227.77"
$ws.Cells.Item(18, 5).Value = "  -1.49%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.33"
$ws.Cells.Item(19, 5).Value = "  -0.92%  "
$ws.Cells.Item(20, 5).Value = "  -2.71%  "
$ws.Cells.Item(21, 5).Value = "  -0.11%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.92"
$ws.Cells.Item(22, 5).Value = "  +0.68%  "
$ws.Cells.Item(23, 5).Value = "  -2.99%  "
$ws.Cells.Item(24, 5).Value = "  -3.30%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "151.25"
$ws.Cells.Item(25, 5).Value = "  -0.48%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "14.75"
$ws.Cells.Item(26, 5).Value = "  -1.78%  "
$ws.Cells.Item(27, 5).Value = "  -1.69%  "
$ws.Cells.Item(28, 5).Value = "  -0.18%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "6.25"
$ws.Cells.Item(29, 5).Value = "  -3.36%  "
$ws.Cells.Item(30, 5).Value = "  -3.28%  "
$ws.Cells.Item(31, 5).Value = "  -4.25%  "
$ws.Cells.Item(32, 5).Value = "  -1.61%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.385.80"
$ws.Cells.Item(33, 5).Value = "  -0.74%  "
$ws.Cells.Item(34, 5).Value = "  -3.11%  "
$ws.Cells.Item(35, 5).Value = "  +2.55%  "
$ws.Cells.Item(37, 5).Value = "  -1.46%  "
$ws.Cells.Item(38, 5).Value = "  -1.26%  "
$ws.Cells.Item(39, 5).Value = "  -2.93%  "
$ws.Cells.Item(40, 5).Value = "  -2.51%  "
$ws.Cells.Item(41, 5).Value = "  +0.79%  "
$ws.Cells.Item(42, 5).Value = "  -0.17%  "
$ws.Cells.Item(43, 5).Value = "  -2.12%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0463"
$ws.Cells.Item(44, 5).Value = "  -1.74%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "5.36"
$ws.Cells.Item(45, 5).Value = "  -2.70%  "
$ws.Cells.Item(46, 5).Value = "  -2.36%  "

# --- Row 47: new entry "RocketPoolETH" inserted; rows 47-50 (old) shift down to 48-51; old row 51 (Cronos) drops off ---
$ws.Cells.Item(47, 2).Value = "RocketPoolETH"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.686.17"
$ws.Cells.Item(47, 5).Value = "  -1.55%  "
$ws.Cells.Item(48, 2).Value = "WEMIXToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.906"
$ws.Cells.Item(48, 5).Value = "  -5.67%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "85.61"
$ws.Cells.Item(49, 5).Value = "  -1.31%  "
$ws.Cells.Item(50, 2).Value = "BitcoinSV"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "42.71"
$ws.Cells.Item(50, 5).Value = "  +7.11%  "
$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0₆0103"
$ws.Cells.Item(51, 5).Value = "  +0.67%  "

Write-Host "cryptos sheet updated"
